# Update cached market-data / profit columns (H-N) on each Leve sheet
# to the latest scheduled-runner snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 5385.5
$ws.Range("I43").Value = 992
$ws.Range("K43").Value = 992
$ws.Range("M43").Value = -923
# Row 62
$ws.Range("H62").Value = 2142
$ws.Range("J62").Value = 2285
$ws.Range("L62").Value = 2285
$ws.Range("N62").Value = -3533
# Row 65
$ws.Range("H65").Value = 2142
$ws.Range("J65").Value = 2285
$ws.Range("L65").Value = 11425
$ws.Range("N65").Value = -17665
# Row 69
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").Value = ""
# Row 72
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").Value = ""

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = ""
# Row 33
$ws.Range("H33").Value = 19950
$ws.Range("I33").Value = 19950
$ws.Range("K33").Value = 19950
$ws.Range("M33").Value = -19621
# Row 36
$ws.Range("H36").Value = 5570.8335
$ws.Range("I36").Value = 5570.8335
$ws.Range("K36").Value = 5570.8335
$ws.Range("M36").Value = -5224.8335
# Row 61
$ws.Range("H61").Value = 1550.7142
$ws.Range("I61").Value = 1550.7142
$ws.Range("K61").Value = 1550.7142
$ws.Range("M61").Value = -1338.7142
# Row 97
$ws.Range("H97").Value = 2008.1052
$ws.Range("I97").Value = 1911.9231
$ws.Range("K97").Value = 1911.9231
$ws.Range("M97").Value = -1415.9231
# Row 108
$ws.Range("H108").Value = 69999
$ws.Range("J108").Value = 69999
$ws.Range("L108").Value = 69999
$ws.Range("N108").Value = -77679
# Row 110
$ws.Range("H110").Value = 1500
$ws.Range("I110").Value = 1500
$ws.Range("K110").Value = 1500
$ws.Range("M110").Value = 545
# Row 136
$ws.Range("H136").Value = 1550.7142
$ws.Range("I136").Value = 1550.7142
$ws.Range("K136").Value = 4652.142599999999
$ws.Range("M136").Value = -2102.142599999999

$ws = $wb.Worksheets.Item("BSM")
# Row 92
$ws.Range("H92").Value = 28000
$ws.Range("J92").Value = 28000
$ws.Range("L92").Value = 28000
$ws.Range("N92").Value = -32992
# Row 94
$ws.Range("H94").Value = 2369.5334
$ws.Range("I94").Value = 1949.3334
$ws.Range("K94").Value = 1949.3334
$ws.Range("M94").Value = -1498.3334
# Row 99
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").Value = ""
# Row 107
$ws.Range("H107").Value = 1111
$ws.Range("J107").Value = 1300
$ws.Range("L107").Value = 1300
$ws.Range("N107").Value = -5140
# Row 109
$ws.Range("H109").Value = 69999
$ws.Range("J109").Value = 69999
$ws.Range("L109").Value = 69999
$ws.Range("N109").Value = -72773

$ws = $wb.Worksheets.Item("CRP")
# Row 41
$ws.Range("H41").Value = 10000
$ws.Range("I41").Value = 10000
$ws.Range("K41").Value = 10000
$ws.Range("M41").Value = -9572
# Row 58
$ws.Range("H58").Value = 607.3333
$ws.Range("I58").Value = 607.3333
$ws.Range("K58").Value = 607.3333
$ws.Range("M58").Value = -404.3333
# Row 109
$ws.Range("H109").Value = 69999
$ws.Range("J109").Value = 69999
$ws.Range("L109").Value = 69999
$ws.Range("N109").Value = -72079
# Row 136
$ws.Range("H136").Value = 607.3333
$ws.Range("I136").Value = 607.3333
$ws.Range("K136").Value = 1821.9999
$ws.Range("M136").Value = 728.0001

$ws = $wb.Worksheets.Item("CUL")
# Row 99
$ws.Range("H99").Value = 5524.8
$ws.Range("I99").Value = 3156
$ws.Range("J99").Value = 15000
$ws.Range("K99").Value = 9468
$ws.Range("L99").Value = 45000
$ws.Range("M99").Value = -7222
$ws.Range("N99").Value = -49492

$ws = $wb.Worksheets.Item("GSM")
# Row 4
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 50
$ws.Range("K4").Value = 50
$ws.Range("M4").Value = 62
# Row 5
$ws.Range("H5").Value = 2000
$ws.Range("I5").Value = 2000
$ws.Range("K5").Value = 2000
$ws.Range("M5").Value = -1888
# Row 9
$ws.Range("H9").Value = 253.5
$ws.Range("I9").Value = 307
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 307
$ws.Range("L9").Value = 200
$ws.Range("M9").Value = -137
$ws.Range("N9").Value = -540
# Row 10
$ws.Range("H10").Value = 500
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 500
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 500
$ws.Range("M10").Value = ""
$ws.Range("N10").Value = -838
# Row 12
$ws.Range("H12").Value = 2003
$ws.Range("I12").Value = 2003
$ws.Range("K12").Value = 2003
$ws.Range("M12").Value = -1863
# Row 14
$ws.Range("H14").Value = 50000000
$ws.Range("I14").Value = 50000000
$ws.Range("K14").Value = 50000000
$ws.Range("M14").Value = -49999832
# Row 19
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").Value = ""
# Row 20
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").Value = ""

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 462041.2
$ws.Range("I40").Value = 4841.4
$ws.Range("J40").Value = 843041
$ws.Range("K40").Value = 4841.4
$ws.Range("L40").Value = 843041
$ws.Range("M40").Value = -4705.4
$ws.Range("N40").Value = -843313
# Row 55
$ws.Range("H55").Value = 1858.4286
$ws.Range("I55").Value = 2092.1667
$ws.Range("K55").Value = 2092.1667
$ws.Range("M55").Value = -1919.1667
# Row 82
$ws.Range("H82").Value = 1299.909
$ws.Range("I82").Value = 1322.1111
$ws.Range("K82").Value = 1322.1111
$ws.Range("M82").Value = -961.1111000000001
# Row 85
$ws.Range("H85").Value = 1299.909
$ws.Range("I85").Value = 1322.1111
$ws.Range("K85").Value = 1322.1111
$ws.Range("M85").Value = -74.11110000000008
# Row 87
$ws.Range("H87").Value = 47999
$ws.Range("J87").Value = 47999
$ws.Range("L87").Value = 47999
$ws.Range("N87").Value = -50245
# Row 90
$ws.Range("H90").Value = 47999
$ws.Range("J90").Value = 47999
$ws.Range("L90").Value = 143997
$ws.Range("N90").Value = -155229
# Row 100
$ws.Range("H100").Value = 2175
$ws.Range("I100").Value = 2175
$ws.Range("K100").Value = 2175
$ws.Range("M100").Value = -1634
# Row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = ""
# Row 132
$ws.Range("H132").Value = 1488.2354
$ws.Range("I132").Value = 1424.7333
$ws.Range("K132").Value = 4274.199900000001
$ws.Range("M132").Value = -1744.199900000001

$ws = $wb.Worksheets.Item("WVR")
# Row 55
$ws.Range("H55").Value = 32292.572
$ws.Range("I55").Value = 12524
$ws.Range("J55").Value = 40200
$ws.Range("K55").Value = 12524
$ws.Range("L55").Value = 40200
$ws.Range("M55").Value = -12247
$ws.Range("N55").Value = -40754
# Row 81
$ws.Range("H81").Value = 3961.875
$ws.Range("I81").Value = 3956.1428
$ws.Range("K81").Value = 7912.2856
$ws.Range("M81").Value = -6851.2856
# Row 84
$ws.Range("H84").Value = 3961.875
$ws.Range("I84").Value = 3956.1428
$ws.Range("K84").Value = 39561.428
$ws.Range("M84").Value = -34257.428
# Row 97
$ws.Range("H97").Value = 23380
$ws.Range("J97").Value = 23380
$ws.Range("L97").Value = 23380
$ws.Range("N97").Value = -25362
# Row 132
$ws.Range("H132").Value = 1034.1666
$ws.Range("I132").Value = 1034.1666
$ws.Range("K132").Value = 3102.4998
$ws.Range("M132").Value = -572.4998000000001
